# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets
# to reflect newly scraped data (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 105
$ws1.Range("F3").Value = 52
$ws1.Range("F6").Value = 9387
$ws1.Range("F7").Value = 842
$ws1.Range("F9").Value = 1193
$ws1.Range("F10").Value = 1130
$ws1.Range("F12").Value = 93
$ws1.Range("F14").Value = 260
$ws1.Range("F15").Value = 415
$ws1.Range("F17").Value = 250
$ws1.Range("F18").Value = 1270

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 105
$ws4.Range("F3").Value = 52
$ws4.Range("F7").Value = 9387
$ws4.Range("F8").Value = 842
$ws4.Range("F10").Value = 1193
$ws4.Range("F11").Value = 1130
$ws4.Range("F13").Value = 93
$ws4.Range("F15").Value = 260
$ws4.Range("F16").Value = 415
$ws4.Range("F18").Value = 250
$ws4.Range("F19").Value = 1270
